$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.704119850187266
$ws1.Range("C2").Value = 0.8205882352941176
$ws1.Range("D2").Value = 0.5224719101123596
$ws1.Range("E2").Value = 0.6384439359267735
$ws1.Range("F2").Value = 0.704119850187266
$ws1.Range("G2").Value = 279
$ws1.Range("H2").Value = 61
$ws1.Range("I2").Value = 473
$ws1.Range("J2").Value = 255

# ---- Sheet 2: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 (label index 15 -> "0.0")
$ws2.Range("B2").Value = 0.6497252747252747
$ws2.Range("C2").Value = 0.8857677902621723
$ws2.Range("D2").Value = 0.7496038034865293

# Row 3 (label index 16 -> "1.0")
$ws2.Range("B3").Value = 0.8205882352941176
$ws2.Range("C3").Value = 0.5224719101123596
$ws2.Range("D3").Value = 0.6384439359267735
$ws2.Range("E3").Value = 534

# Row 4 (label index 17 -> "accuracy")
$ws2.Range("B4").Value = 0.704119850187266
$ws2.Range("C4").Value = 0.704119850187266
$ws2.Range("D4").Value = 0.704119850187266
$ws2.Range("E4").Value = 0.704119850187266

# Row 5 (label index 18 -> "macro avg")
$ws2.Range("B5").Value = 0.7351567550096962
$ws2.Range("C5").Value = 0.704119850187266
$ws2.Range("D5").Value = 0.6940238697066514
$ws2.Range("E5").Value = 1068

# Row 6 (label index 19 -> "weighted avg")
$ws2.Range("B6").Value = 0.7351567550096962
$ws2.Range("C6").Value = 0.704119850187266
$ws2.Range("D6").Value = 0.6940238697066513
$ws2.Range("E6").Value = 1068

# ---- Sheet 3: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 473
$ws3.Range("C2").Value = 61
$ws3.Range("B3").Value = 255
$ws3.Range("C3").Value = 279
